$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "66.039.87"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +0.19%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.329.39"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +1.56%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Range("E4").Value = "  -0.15%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "186.72"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.91%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "558.43"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.21%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -0.06%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "3.323.57"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +1.51%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.578"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -2.55%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.177"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -5.55%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.580"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -1.33%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "46.07"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -2.74%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000264"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -0.99%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "3.846.19"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +1.30%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "8.48"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -1.31%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "574.89"
$cell.ClearFormats()
$ws.Range("E16").Value = "  -8.97%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "66.027.39"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  +0.58%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "3.310.46"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.42%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.76"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -0.78%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "10.91"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -4.06%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.894"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -1.14%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "18.11"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +1.74%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "5.00"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +1.11%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "98.10"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -8.38%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "3.97"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -0.33%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.69"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +1.06%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "9.42"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -1.18%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "8.46"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -2.77%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "30.64"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +0.64%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "6.73"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +6.94%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.72"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -8.24%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "566.26"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +4.50%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "10.84"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("E35").Value = "  -1.45%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "3.722.53"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +0.10%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "55.61"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -2.89%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "34.15"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +4.68%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.127"
$cell.ClearFormats()
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0691"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -5.49%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.62"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "3.13"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -8.61%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "3.36"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +3.18%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.335"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.73%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0408"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("B48").Value = "CoreDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.93"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -12.06%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +0.14%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.52"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -4.03%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "126.92"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +3.96%  "
